{"js": "// Update the date line and all 26 \"three-digit / one-digit\" division\n// answers in the table to the new values from the latest generated output.\n\nconst replacements = [\n    { old: \"2024-09-15 Sunday\", new: \"2024-09-16 Monday\" },\n    { old: \"221\u00f74=55, 1\",       new: \"127\u00f77=18, 1\" },\n    { old: \"122\u00f72=61, 0\",       new: \"437\u00f77=62, 3\" },\n    { old: \"586\u00f79=65, 1\",       new: \"335\u00f79=37, 2\" },\n    { old: \"151\u00f74=37, 3\",       new: \"808\u00f77=115, 3\" },\n    { old: \"481\u00f76=80, 1\",       new: \"734\u00f76=122, 2\" },\n    { old: \"535\u00f72=267, 1\",      new: \"469\u00f74=117, 1\" },\n    { old: \"211\u00f76=35, 1\",       new: \"784\u00f75=156, 4\" },\n    { old: \"324\u00f72=162, 0\",      new: \"208\u00f75=41, 3\" },\n    { old: \"410\u00f74=102, 2\",      new: \"834\u00f77=119, 1\" },\n    { old: \"245\u00f76=40, 5\",       new: \"697\u00f77=99, 4\" },\n    { old: \"843\u00f78=105, 3\",      new: \"897\u00f74=224, 1\" },\n    { old: \"702\u00f74=175, 2\",      new: \"767\u00f74=191, 3\" },\n    { old: \"103\u00f72=51, 1\",       new: \"113\u00f79=12, 5\" },\n    { old: \"853\u00f72=426, 1\",      new: \"380\u00f78=47, 4\" },\n    { old: \"529\u00f77=75, 4\",       new: \"498\u00f77=71, 1\" },\n    { old: \"127\u00f75=25, 2\",       new: \"888\u00f76=148, 0\" },\n    { old: \"718\u00f79=79, 7\",       new: \"489\u00f72=244, 1\" },\n    { old: \"255\u00f78=31, 7\",       new: \"302\u00f74=75, 2\" },\n    { old: \"165\u00f76=27, 3\",       new: \"202\u00f74=50, 2\" },\n    { old: \"812\u00f77=116, 0\",      new: \"924\u00f75=184, 4\" },\n    { old: \"468\u00f72=234, 0\",      new: \"120\u00f79=13, 3\" },\n    { old: \"488\u00f73=162, 2\",      new: \"881\u00f76=146, 5\" },\n    { old: \"278\u00f79=30, 8\",       new: \"893\u00f72=446, 1\" },\n    { old: \"128\u00f79=14, 2\",       new: \"677\u00f78=84, 5\" },\n    { old: \"819\u00f74=204, 3\",      new: \"498\u00f76=83, 0\" },\n];\n\nconst body = context.document.body;\n\nfor (const pair of replacements) {\n    const results = body.search(pair.old, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(pair.new, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Update the date line and all 26 \"three-digit / one-digit\" division\n# answers in the table to the new values from the latest generated output.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = '2024-09-15 Sunday'; new = '2024-09-16 Monday'},\n    @{old = '221\u00f74=55, 1';       new = '127\u00f77=18, 1'},\n    @{old = '122\u00f72=61, 0';       new = '437\u00f77=62, 3'},\n    @{old = '586\u00f79=65, 1';       new = '335\u00f79=37, 2'},\n    @{old = '151\u00f74=37, 3';       new = '808\u00f77=115, 3'},\n    @{old = '481\u00f76=80, 1';       new = '734\u00f76=122, 2'},\n    @{old = '535\u00f72=267, 1';      new = '469\u00f74=117, 1'},\n    @{old = '211\u00f76=35, 1';       new = '784\u00f75=156, 4'},\n    @{old = '324\u00f72=162, 0';      new = '208\u00f75=41, 3'},\n    @{old = '410\u00f74=102, 2';      new = '834\u00f77=119, 1'},\n    @{old = '245\u00f76=40, 5';       new = '697\u00f77=99, 4'},\n    @{old = '843\u00f78=105, 3';      new = '897\u00f74=224, 1'},\n    @{old = '702\u00f74=175, 2';      new = '767\u00f74=191, 3'},\n    @{old = '103\u00f72=51, 1';       new = '113\u00f79=12, 5'},\n    @{old = '853\u00f72=426, 1';      new = '380\u00f78=47, 4'},\n    @{old = '529\u00f77=75, 4';       new = '498\u00f77=71, 1'},\n    @{old = '127\u00f75=25, 2';       new = '888\u00f76=148, 0'},\n    @{old = '718\u00f79=79, 7';       new = '489\u00f72=244, 1'},\n    @{old = '255\u00f78=31, 7';       new = '302\u00f74=75, 2'},\n    @{old = '165\u00f76=27, 3';       new = '202\u00f74=50, 2'},\n    @{old = '812\u00f77=116, 0';      new = '924\u00f75=184, 4'},\n    @{old = '468\u00f72=234, 0';      new = '120\u00f79=13, 3'},\n    @{old = '488\u00f73=162, 2';      new = '881\u00f76=146, 5'},\n    @{old = '278\u00f79=30, 8';       new = '893\u00f72=446, 1'},\n    @{old = '128\u00f79=14, 2';       new = '677\u00f78=84, 5'},\n    @{old = '819\u00f74=204, 3';      new = '498\u00f76=83, 0'}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.old\n    $find.Replacement.Text = $pair.new\n    $find.Execute($pair.old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.new, 2)\n}\n"}
